$wb = $excel.ActiveWorkbook

# "Generate Report for Handback" — refresh the timestamps / priority that
# were recorded for the 88a8d918... handback file, and for the shared
# string entries that f6b21a53... (row 5) happens to reuse.

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-31 04:15:48"
$wsOverview.Range("G5").Value = "2016-08-31 04:15:48"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H3").Value = "2016-08-31 04:15:42"
$wsZhCn.Range("K3").Value = "2016-08-31 04:16:16"
$wsZhCn.Range("E5").Value = "mt"
$wsZhCn.Range("H5").Value = "2016-08-31 04:15:42"
$wsZhCn.Range("K5").Value = "2016-08-31 04:16:16"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("H3").Value = "2016-08-31 04:15:48"
$wsDeDe.Range("K3").Value = "2016-08-31 04:16:23"
$wsDeDe.Range("E5").Value = "mt"
$wsDeDe.Range("H5").Value = "2016-08-31 04:15:48"
$wsDeDe.Range("K5").Value = "2016-08-31 04:16:23"
